$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.860.02'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '1.631.18'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  +0.69%  '

$ws.Range("D5").Value = '''214.34'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").Value = '''0.0632'
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").Value = '''19.54'
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").Value = '''0.0792'
$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '1.857.42'
$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.25'
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.624.94'
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("D15").Value = '''0.544'
$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").Value = '0.0₃0754'
$ws.Range("E16").Value = '  -0.57%  '

$ws.Range("D17").Value = '''62.62'
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").Value = '25.874.82'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("E20").Value = '  -0.95%  '

$ws.Range("D21").Value = '''193.02'
$ws.Range("E21").Value = '  +0.97%  '

$ws.Range("D22").Value = '''9.91'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").Value = '''6.24'
$ws.Range("E23").Value = '  -0.41%  '

$ws.Range("D24").Value = '''1.82'
$ws.Range("E24").Value = '  +1.27%  '

$ws.Range("D25").Value = '''143.29'
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("E26").Value = '  +0.44%  '

$ws.Range("E27").Value = '  +2.91%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").Value = '''15.42'
$ws.Range("E29").Value = '  -0.78%  '

$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("D31").Value = '''0.0497'
$ws.Range("E31").Value = '  +0.80%  '

$ws.Range("E32").Value = '  -0.75%  '

$ws.Range("D33").Value = '''3.22'
$ws.Range("E33").Value = '  +0.19%  '

$ws.Range("E34").Value = '  -1.85%  '

$ws.Range("D35").Value = '''2.43'
$ws.Range("E35").Value = '  +1.66%  '

$ws.Range("D36").Value = '''0.901'
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").Value = '1.137.98'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").Value = '''0.546'
$ws.Range("E38").Value = '  +0.42%  '

$ws.Range("E39").Value = '  -1.20%  '

$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''99.08'
$ws.Range("E42").Value = '  -1.36%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.44'
$ws.Range("E43").Value = '  -2.47%  '

$ws.Range("D44").Value = '''0.794'
$ws.Range("E44").Value = '  -0.56%  '

$ws.Range("D45").Value = '1.767.24'
$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("E46").Value = '  +3.00%  '

$ws.Range("D47").Value = '''56.24'
$ws.Range("E47").Value = '  +1.53%  '

$ws.Range("D48").Value = '''0.0528'
$ws.Range("E48").Value = '  +3.12%  '

$ws.Range("E49").Value = '  -1.39%  '

$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").Value = '''7.62'
$ws.Range("E51").Value = '  +0.52%  '
